# Update the "想去人数" (number of people interested) counts that changed
# between two generated outputs, as reflected in the gh-pages data refresh.
#
# Sheet "展览" (Exhibition, sheet1): F3 825->826, F5 1016->1017, F6 2378->2379
# Sheet "全部类型" (All Types, sheet4): F3 825->826, F7 1016->1017, F8 2378->2379

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value = 826
$wsExhibition.Range("F5").Value = 1017
$wsExhibition.Range("F6").Value = 2379

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 826
$wsAll.Range("F7").Value = 1017
$wsAll.Range("F8").Value = 2379
